$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 first (most recently added supplement, 2011) - filled left to right
$ws.Range("A10").Value = 2011
$ws.Range("B10").Value = "ブルーフォレスト通信3"
$ws.Range("C10").Value = "Blue Forest Report 3"
$ws.Range("D10").Value = "Grimpeur"
$ws.Range("E10").Value = "blue_forest_report3.jpg"
$ws.Range("F10").Value = "supplement"

# Row 8 / row 9 (2010 supplements)
$ws.Range("A8").Value = 2010
$ws.Range("A9").Value = 2010

$ws.Range("E8").Value = "blue_forest_report.jpg"
$ws.Range("E9").Value = "blue_forest_report2.jpg"

$ws.Range("B9").Value = "ブルーフォレスト通信2"
$ws.Range("B8").Value = "ブルーフォレスト通信"

$ws.Range("C8").Value = "Blue Forest Report"
$ws.Range("C9").Value = "Blue Forest Report 2"

$ws.Range("D8").Value = "Grimpeur"
$ws.Range("D9").Value = "Grimpeur"

$ws.Range("F8").Value = "supplement"
$ws.Range("F9").Value = "supplement"

# Row 5 grew by one point once the new rows were in place (content reflow)
$ws.Rows.Item(5).RowHeight = 17

# Selection ends up on the newly-filled image column for the added rows
$ws.Range("F8:F10").Select() | Out-Null
